# Add the "Week 12" sheet: a new week of games was added to the front of
# the workbook (matching the existing "Week N" tabs pattern), pushing the
# previous sheets one slot to the right. Excel inserts a new sheet before
# the currently-active sheet, which was the first tab ("Week 11"), so the
# new sheet lands in the very first position - exactly where "Week 12"
# needs to be.

$wb = $excel.ActiveWorkbook
$week12 = $wb.Worksheets.Add()
$week12.Name = "Week 12"

# Same 3-column layout (game / temp / wind) used by every other weekly
# sheet in this workbook.
$data = @(
    @("game", "temp", "wind"),
    @("BUFvsHOU", 75, 8),
    @("CHIvsPIT", 49, 7),
    @("CINvsNE", 55, 6),
    @("DETvsNYG", 49, 7),
    @("GBvsMIN", 46, 8),
    @("SEAvsTEN", 61, 5),
    @("INDvsKC", 57, 5),
    @("BALvsNYJ", 54, 7),
    @("CLEvsLV", 61, 6),
    @("ARIvsJAX", 65, 3),
    @("DALvsPHI", 68, 6),
    @("ATLvsNO", 74, 5),
    @("LAvsTB", 59, 8),
    @("CARvsSF", 57, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $week12.Range("A$row").Value = $data[$i][0]
    $week12.Range("B$row").Value = $data[$i][1]
    $week12.Range("C$row").Value = $data[$i][2]
}

# Leave the selection where the prior sheet's last save left it, matching
# the carried-over cursor position.
[void]$week12.Range("C16").Select()
